$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 22222444
$ws.Range("I11").Value = 22222444
$ws.Range("K11").Value = 22222444
$ws.Range("M11").Value = -22222304

$ws.Range("H28").Value = 2116045.8
$ws.Range("I28").Value = 3800184.5
$ws.Range("J28").Value = 10872.25
$ws.Range("K28").Value = 3800184.5
$ws.Range("L28").Value = 10872.25
$ws.Range("M28").Value = -3799699.5
$ws.Range("N28").Value = -11842.25

$ws.Range("H31").Value = 400
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H69").Value = 2866.6667
$ws.Range("I69").Value = 2866.6667
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 8600.000100000001
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -7726.000100000001
$ws.Range("N69").ClearContents()

$ws.Range("H70").Value = 12228
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 13569
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 40707
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -41247

$ws.Range("H72").Value = 2866.6667
$ws.Range("I72").Value = 2866.6667
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 25800.0003
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -21432.0003
$ws.Range("N72").ClearContents()

$ws.Range("H73").Value = 12228
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 13569
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 40707
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -42579

$ws.Range("H137").Value = 33745.613
$ws.Range("I137").Value = 1409.174
$ws.Range("J137").Value = 126712.875
$ws.Range("K137").Value = 4227.522
$ws.Range("L137").Value = 380138.625
$ws.Range("M137").Value = -1677.522
$ws.Range("N137").Value = -385238.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2601.25
$ws.Range("I61").Value = 1326.3572
$ws.Range("K61").Value = 1326.3572
$ws.Range("M61").Value = -1114.3572

$ws.Range("H74").Value = 1249.919
$ws.Range("I74").Value = 1023.96295
$ws.Range("K74").Value = 1023.96295
$ws.Range("M74").Value = -149.96295

$ws.Range("H77").Value = 1249.919
$ws.Range("I77").Value = 1023.96295
$ws.Range("K77").Value = 5119.81475
$ws.Range("M77").Value = -751.8147499999995

$ws.Range("H132").Value = 1956.1613
$ws.Range("J132").Value = 2522.2307
$ws.Range("L132").Value = 7566.6921
$ws.Range("N132").Value = -12626.6921

$ws.Range("H136").Value = 2601.25
$ws.Range("I136").Value = 1326.3572
$ws.Range("K136").Value = 3979.0716
$ws.Range("M136").Value = -1429.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H134").Value = 7351.625
$ws.Range("I134").Value = 8785.526
$ws.Range("K134").Value = 26356.578
$ws.Range("M134").Value = -23821.578

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1991.1666
$ws.Range("I31").Value = 1361
$ws.Range("K31").Value = 1361
$ws.Range("M31").Value = -1066

$ws.Range("H34").Value = 1991.1666
$ws.Range("I34").Value = 1361
$ws.Range("K34").Value = 1361
$ws.Range("M34").Value = -1159

$ws.Range("H132").Value = 2737.25
$ws.Range("I132").Value = 1444.4166
$ws.Range("K132").Value = 4333.2498
$ws.Range("M132").Value = -1803.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 551.8889
$ws.Range("I26").Value = 543.6667
$ws.Range("J26").Value = 568.3333
$ws.Range("K26").Value = 1631.0001
$ws.Range("L26").Value = 1704.9999
$ws.Range("M26").Value = -1343.0001
$ws.Range("N26").Value = -2280.9999

$ws.Range("H36").Value = 709.8
$ws.Range("I36").Value = 1100
$ws.Range("J36").Value = 124.5
$ws.Range("K36").Value = 3300
$ws.Range("L36").Value = 373.5
$ws.Range("M36").Value = -3131
$ws.Range("N36").Value = -711.5

$ws.Range("H138").Value = 2931.4
$ws.Range("J138").Value = 5275
$ws.Range("L138").Value = 15825
$ws.Range("N138").Value = -26105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13188.777
$ws.Range("I70").Value = 20039.8
$ws.Range("K70").Value = 20039.8
$ws.Range("M70").Value = -19769.8

$ws.Range("H73").Value = 13188.777
$ws.Range("I73").Value = 20039.8
$ws.Range("K73").Value = 20039.8
$ws.Range("M73").Value = -19103.8

$ws.Range("H80").Value = 3891
$ws.Range("J80").Value = 3933
$ws.Range("L80").Value = 3933
$ws.Range("N80").Value = -5929

$ws.Range("H83").Value = 3891
$ws.Range("J83").Value = 3933
$ws.Range("L83").Value = 19665
$ws.Range("N83").Value = -29649

$ws.Range("H113").Value = 1398.5
$ws.Range("I113").Value = 1297
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1297
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 873
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 1242537.9
$ws.Range("I132").Value = 1749582
$ws.Range("J132").Value = 3096.4443
$ws.Range("K132").Value = 5248746
$ws.Range("L132").Value = 9289.332900000001
$ws.Range("M132").Value = -5246216
$ws.Range("N132").Value = -14349.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1909.9688
$ws.Range("I132").Value = 1933.6
$ws.Range("K132").Value = 5800.799999999999
$ws.Range("M132").Value = -3270.799999999999

$ws.Range("H136").Value = 2343.4644
$ws.Range("I136").Value = 1386.5238
$ws.Range("K136").Value = 4159.5714
$ws.Range("M136").Value = -1609.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 640.03845
$ws.Range("I107").Value = 560.6818
$ws.Range("K107").Value = 1682.0454
$ws.Range("M107").Value = 237.9546

$ws.Range("H122").Value = 61446
$ws.Range("I122").Value = 98321.75
$ws.Range("J122").Value = 2444.8
$ws.Range("K122").Value = 294965.25
$ws.Range("L122").Value = 7334.400000000001
$ws.Range("M122").Value = -292515.25
$ws.Range("N122").Value = -12234.4

$ws.Range("H132").Value = 1095.8788
$ws.Range("I132").Value = 867.7586
$ws.Range("J132").Value = 2749.75
$ws.Range("K132").Value = 2603.2758
$ws.Range("L132").Value = 8249.25
$ws.Range("M132").Value = -73.27579999999989
$ws.Range("N132").Value = -13309.25

$ws.Range("H136").Value = 22224928
$ws.Range("I136").Value = 29242458
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 87727374
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -87724824
$ws.Range("N136").Value = -13350
